$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1901
$ws1.Range("F4").Value = 866
$ws1.Range("F6").Value = 51
$ws1.Range("G6").Value = 50
$ws1.Range("F10").Value = 158
$ws1.Range("F13").Value = 4493
$ws1.Range("F16").Value = 492
$ws1.Range("F20").Value = 1218
$ws1.Range("F21").Value = 2277
$ws1.Range("F23").Value = 63
$ws1.Range("F24").Value = 42
$ws1.Range("F25").Value = 54
$ws1.Range("F26").Value = 2200
$ws1.Range("F29").Value = 21
$ws1.Range("F30").Value = 157
$ws1.Range("F34").Value = 34

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1901
$ws4.Range("F4").Value = 866
$ws4.Range("F6").Value = 51
$ws4.Range("G6").Value = 50
$ws4.Range("F10").Value = 158
$ws4.Range("F14").Value = 4493
$ws4.Range("F17").Value = 492
$ws4.Range("F21").Value = 1218
$ws4.Range("F22").Value = 2277
$ws4.Range("F24").Value = 63
$ws4.Range("F25").Value = 42
$ws4.Range("F26").Value = 54
$ws4.Range("G26").Value = 50
$ws4.Range("F27").Value = 2200
$ws4.Range("F30").Value = 21
$ws4.Range("F31").Value = 157
$ws4.Range("F35").Value = 34
